$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (column C) date for every data row (2..517) from 45203 to 45204
$ws.Range("C2:C517").Value = 45204

# 2. Give row 517 an explicit row height (matches the diff: ht="15" customHeight="1")
$ws.Rows.Item(517).RowHeight = 15

# 3. Append a new row (518) with the new logging notification
$ws.Cells.Item(518, 1).Value = "A 47589-2023"

$ws.Cells.Item(518, 2).Value = 45203
$ws.Cells.Item(518, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(518, 3).Value = 45204
$ws.Cells.Item(518, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(518, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(518, 5).Value = "FALKENBERG"

$ws.Cells.Item(518, 7).Value = 0.6
$ws.Cells.Item(518, 8).Value = 0
$ws.Cells.Item(518, 9).Value = 0
$ws.Cells.Item(518, 10).Value = 0
$ws.Cells.Item(518, 11).Value = 0
$ws.Cells.Item(518, 12).Value = 0
$ws.Cells.Item(518, 13).Value = 0
$ws.Cells.Item(518, 14).Value = 0
$ws.Cells.Item(518, 15).Value = 0
$ws.Cells.Item(518, 16).Value = 0
$ws.Cells.Item(518, 17).Value = 0

# R518 keeps the same wrap-text styling used throughout column R
$ws.Cells.Item(518, 18).WrapText = $true
